# Add new columns I ("I0") and J ("IF") to the worksheet, mirroring the
# style of the existing header cells (e.g. H1) and filling in the data
# values for rows 2 through 41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header styling used by the existing headers (e.g. column H) by
# copying its format, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-41 (row, I-value, J-value)
$data = @(
    @(2, 7, 8),
    @(3, 5, 6),
    @(4, 6, 6),
    @(5, 6, 7),
    @(6, 7, 7),
    @(7, 6, 6),
    @(8, 8, 9),
    @(9, 6, 6),
    @(10, 8, 8),
    @(11, 7, 7),
    @(12, 6, 7),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 1, 2),
    @(18, 5, 6),
    @(19, 5, 6),
    @(20, 6, 6),
    @(21, 7, 7),
    @(22, 7, 8),
    @(23, 6, 7),
    @(24, 8, 8),
    @(25, 9, 9),
    @(26, 7, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 9, 9),
    @(31, 6, 6),
    @(32, 9, 9),
    @(33, 5, 5),
    @(34, 6, 8),
    @(35, 7, 7),
    @(36, 6, 6),
    @(37, 4, 5),
    @(38, 4, 4),
    @(39, 4, 4),
    @(40, 5, 5),
    @(41, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
